$d = $word.ActiveDocument

function Set-ParagraphRunXml($Paragraph, $InnerRunXml) {
    $pRange = $Paragraph.Range
    # Exclude the trailing paragraph mark (End - 1) so the paragraph
    # itself (its identity / pPr) is left untouched - only its runs change.
    $contentRange = $d.Range($pRange.Start, $pRange.End - 1)

    $package = @"
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>$InnerRunXml</w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

    $contentRange.InsertXML($package)
}

# --- Change 1: paragraph 1 - merge "D" + "efault value" + ", skill points,
# class modifier" runs (after the line break) into a single run, leaving
# "Health modifiers" as its own run. ---
$para1 = $d.Paragraphs.Item(1)
$run1Xml = '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Health modifiers</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:br/><w:t>Default value, skill points, class modifier</w:t></w:r>'
Set-ParagraphRunXml $para1 $run1Xml

# --- Change 2: paragraph 2 - split the single run into several runs and
# change the default health value from 500-1000 to 2000-10000 (and the
# per-level bonus from 100 to 1000). ---
$para2 = $d.Paragraphs.Item(2)
$run2Xml = '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Default value: </w:t></w:r>' + `
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>2000</w:t></w:r>' + `
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>-1000</w:t></w:r>' + `
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>0</w:t></w:r>' + `
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> heath + 100</w:t></w:r>' + `
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>0</w:t></w:r>' + `
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> health for every player level</w:t></w:r>'
Set-ParagraphRunXml $para2 $run2Xml

# --- Change 3: paragraph 3 - bump the skill points base value from 75 to
# 100, keeping every other run exactly as it was. ---
$para3 = $d.Paragraphs.Item(3)
$run3Xml = '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Skill points: </w:t></w:r>' + `
           '<w:r w:rsidR="00AF59E9"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>100</w:t></w:r>' + `
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> health </w:t></w:r>' + `
           '<w:r w:rsidR="00AF59E9"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">* skill points invested </w:t></w:r>' + `
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">+ </w:t></w:r>' + `
           '<w:r w:rsidR="00AF59E9"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>200</w:t></w:r>' + `
           '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>% of number of skill points invested in health raised by 2</w:t></w:r>'
Set-ParagraphRunXml $para3 $run3Xml
